$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 7, shifting existing rows 7-25 down to 8-26.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with its data.
$ws.Range("A7").Value = 9
$ws.Range("B7").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C7").Value = "Metropolitana"
$ws.Range("D7").Value = 44453
$ws.Range("D7").NumberFormat = $ws.Range("D8").NumberFormat
$ws.Range("E7").Value = 13
$ws.Range("F7").Value = 100112035
$ws.Range("G7").Value = "Bruselas (repollito)"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 25
$ws.Range("K7").Value = 25000
$ws.Range("L7").Value = 26000
$ws.Range("M7").Value = 25520
$ws.Range("N7").Value = "$/malla 15 kilos"
$ws.Range("O7").Value = "Hijuelas"
$ws.Range("P7").Value = 1701
$ws.Range("Q7").Value = 15
$ws.Range("R7").Value = "Hortaliza"
